$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new order rows (13 and 14) below the existing data, mirroring
# the "Pending" / no-delivery-assigned-yet layout already used by row 12.

$ws.Range("A13").Value = "89bdc2f6-0e22-47a8-b4f2-b7b5696fc495"
$ws.Range("B13").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C13").Value = "P2001"
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "Pending"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "2025-08-10"
$ws.Range("G13").ClearFormats()

$ws.Range("A14").Value = "0947da20-6ab3-444d-97b4-2aa9c1662a75"
$ws.Range("B14").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("C14").Value = "P2005"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = "Pending"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "2025-08-10"
$ws.Range("G14").ClearFormats()
